$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1919.375
$ws.Range("I6").Value = 88.75
$ws.Range("J6").Value = 3750
$ws.Range("K6").Value = 266.25
$ws.Range("L6").Value = 11250
$ws.Range("M6").Value = -154.25
$ws.Range("N6").Value = -11474
# Row 40
$ws.Range("H40").Value = 1200
$ws.Range("J40").Value = 1400
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1750
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
# Row 132
$ws.Range("H132").Value = 91830860
$ws.Range("I132").Value = 111124940
$ws.Range("K132").Value = 333374820
$ws.Range("M132").Value = -333372290
# Row 138
$ws.Range("H138").Value = 3985.182
$ws.Range("I138").Value = 805.087
$ws.Range("J138").Value = 4947.579
$ws.Range("K138").Value = 2415.261
$ws.Range("L138").Value = 14842.737
$ws.Range("M138").Value = 2724.739
$ws.Range("N138").Value = -25122.737
# Row 141
$ws.Range("H141").Value = 9218.481
$ws.Range("I141").Value = 10809.333
$ws.Range("J141").Value = 3650.5
$ws.Range("K141").Value = 32427.999
$ws.Range("L141").Value = 10951.5
$ws.Range("M141").Value = -27247.999
$ws.Range("N141").Value = -21311.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 14893.667
$ws.Range("J6").Value = 14893.667
$ws.Range("L6").Value = 14893.667
$ws.Range("N6").Value = -15239.667
# Row 32
$ws.Range("H32").Value = 6004.317
$ws.Range("I32").Value = 5247.6484
$ws.Range("J32").Value = 13003.5
$ws.Range("K32").Value = 5247.6484
$ws.Range("L32").Value = 13003.5
$ws.Range("M32").Value = -4960.6484
$ws.Range("N32").Value = -13577.5
# Row 61
$ws.Range("H61").Value = 1093.3334
$ws.Range("J61").Value = 1429.125
$ws.Range("L61").Value = 1429.125
$ws.Range("N61").Value = -1853.125
# Row 136
$ws.Range("H136").Value = 1093.3334
$ws.Range("J136").Value = 1429.125
$ws.Range("L136").Value = 4287.375
$ws.Range("N136").Value = -9387.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 178.96666
$ws.Range("I80").Value = 63.4
$ws.Range("J80").Value = 236.75
$ws.Range("K80").Value = 63.4
$ws.Range("L80").Value = 236.75
$ws.Range("M80").Value = 934.6
$ws.Range("N80").Value = -2232.75
# Row 83
$ws.Range("H83").Value = 178.96666
$ws.Range("I83").Value = 63.4
$ws.Range("J83").Value = 236.75
$ws.Range("K83").Value = 317
$ws.Range("L83").Value = 1183.75
$ws.Range("M83").Value = 4675
$ws.Range("N83").Value = -11167.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 11116177
$ws.Range("I99").Value = 25003124
$ws.Range("J99").Value = 6619
$ws.Range("K99").Value = 25003124
$ws.Range("L99").Value = 6619
$ws.Range("M99").Value = -25001626
$ws.Range("N99").Value = -9615
# Row 126
$ws.Range("H126").Value = 11116177
$ws.Range("I126").Value = 25003124
$ws.Range("J126").Value = 6619
$ws.Range("K126").Value = 75009372
$ws.Range("L126").Value = 19857
$ws.Range("M126").Value = -75006902
$ws.Range("N126").Value = -24797

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 23256.857
$ws.Range("I11").Value = 24159.6
$ws.Range("J11").Value = 21000
$ws.Range("K11").Value = 72478.79999999999
$ws.Range("L11").Value = 63000
$ws.Range("M11").Value = -72338.79999999999
$ws.Range("N11").Value = -63280
# Row 12
$ws.Range("H12").Value = 132.41667
$ws.Range("I12").Value = 13.333333
$ws.Range("J12").Value = 172.11111
$ws.Range("K12").Value = 39.999999
$ws.Range("L12").Value = 516.3333299999999
$ws.Range("M12").Value = 133.000001
$ws.Range("N12").Value = -862.3333299999999
# Row 92
$ws.Range("H92").Value = 281.875
$ws.Range("I92").Value = 287.85715
$ws.Range("J92").Value = 240
$ws.Range("K92").Value = 863.5714499999999
$ws.Range("L92").Value = 720
$ws.Range("M92").Value = 384.4285500000001
$ws.Range("N92").Value = -3216
# Row 131
$ws.Range("H131").Value = 790.73193
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 804.5268600000001
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 2413.58058
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -12493.58058
# Row 136
$ws.Range("H136").Value = 2835.0386
$ws.Range("I136").Value = 2431.9333
$ws.Range("J136").Value = 3384.7273
$ws.Range("K136").Value = 7295.7999
$ws.Range("L136").Value = 10154.1819
$ws.Range("M136").Value = -2195.7999
$ws.Range("N136").Value = -20354.1819
# Row 137
$ws.Range("H137").Value = 1144.6842
$ws.Range("I137").Value = 857.0714
$ws.Range("J137").Value = 1950
$ws.Range("K137").Value = 2571.2142
$ws.Range("L137").Value = 5850
$ws.Range("M137").Value = 2528.7858
$ws.Range("N137").Value = -16050

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6022
$ws.Range("I70").Value = 5555.5864
$ws.Range("K70").Value = 5555.5864
$ws.Range("M70").Value = -5285.5864
# Row 73
$ws.Range("H73").Value = 6022
$ws.Range("I73").Value = 5555.5864
$ws.Range("K73").Value = 5555.5864
$ws.Range("M73").Value = -4619.5864
# Row 80
$ws.Range("H80").Value = 22729628
$ws.Range("I80").Value = 41668820
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 41668820
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -41667822
$ws.Range("N80").Value = -4596
# Row 83
$ws.Range("H83").Value = 22729628
$ws.Range("I83").Value = 41668820
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 208344100
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -208339108
$ws.Range("N83").Value = -22984

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20224
# Row 15
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20340
# Row 16
$ws.Range("H16").Value = 1254.2858
$ws.Range("I16").Value = 1120
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1120
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -950
$ws.Range("N16").Value = -3340
# Row 40
$ws.Range("H40").Value = 8395.462
$ws.Range("I40").Value = 8066.3
$ws.Range("J40").Value = 9492.666999999999
$ws.Range("K40").Value = 8066.3
$ws.Range("L40").Value = 9492.666999999999
$ws.Range("M40").Value = -7930.3
$ws.Range("N40").Value = -9764.666999999999
# Row 68
$ws.Range("H68").Value = 650.5700000000001
$ws.Range("I68").Value = 650.5700000000001
$ws.Range("K68").Value = 650.5700000000001
$ws.Range("M68").Value = 98.42999999999995
# Row 71
$ws.Range("H71").Value = 650.5700000000001
$ws.Range("I71").Value = 650.5700000000001
$ws.Range("K71").Value = 3252.85
$ws.Range("M71").Value = 491.1499999999996

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 33934.5
$ws.Range("J47").Value = 33934.5
$ws.Range("L47").Value = 33934.5
$ws.Range("N47").Value = -35078.5
# Row 100
$ws.Range("H100").Value = 472.125
$ws.Range("I100").Value = 468.14285
$ws.Range("K100").Value = 936.2857
$ws.Range("M100").Value = -395.2857
# Row 133
$ws.Range("H133").Value = 49330
$ws.Range("J133").Value = 49330
$ws.Range("L133").Value = 49330
$ws.Range("N133").Value = -59450
